# Generated Reports to be added
# Adds a "Clone Detection Tools" report table (D24:J30) below the existing
# summary table, bolds the chart title text "Average Detection Accuracy",
# and sets the page orientation to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlLeft   = -4131

# ---------------------------------------------------------------------
# 1. New header row (row 24) - "Clone Detection Tools" report header
# ---------------------------------------------------------------------
$ws.Range("D24").Value = "Clone Detection Tools"
$ws.Range("E24").Value = "Changes at the same time"
$ws.Range("F24").Value = "Cloned Cochnage"
$ws.Range("G24").Value = "Detected Cochange by Clone Detection Tools"
$ws.Range("J24").Value = "Detection Accuracy"

$headerRange = $ws.Range("D24:J24")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = $xlCenter
$headerRange.VerticalAlignment = $xlCenter
$headerRange.WrapText = $true

$ws.Range("G24:I24").Merge()
$ws.Rows.Item(24).RowHeight = 60

# ---------------------------------------------------------------------
# 2. Data rows (25-30) - one row per clone-detection tool
# ---------------------------------------------------------------------
$ws.Range("D25").Value = "deckard"
$ws.Range("D26").Value = "nicad5"
$ws.Range("D27").Value = "conqat"
$ws.Range("D28").Value = "simcad"
$ws.Range("D29").Value = "iclones"
$ws.Range("D30").Value = "simian"

$dRange = $ws.Range("D25:D30")
$dRange.Font.Bold = $true
$dRange.HorizontalAlignment = $xlCenter
$dRange.VerticalAlignment = $xlCenter

# "Changes at the same time" / "Cloned Cochnage" columns (merged vertically)
$ws.Range("E25").Value = "c1 to c20"
$ws.Range("F25").Value = "c1 to c10"

$efRange = $ws.Range("E25:F30")
$efRange.HorizontalAlignment = $xlCenter
$efRange.VerticalAlignment = $xlCenter
$efRange.WrapText = $true

$ws.Range("E25:E30").Merge()
$ws.Range("F25:F30").Merge()

# "Detected Cochange by Clone Detection Tools" column (merged G:I per row)
$ws.Range("G25").Value = "c1, c3, c5, c6"
$ws.Range("G26").Value = "c2, c3, c4, c6, c7"
$ws.Range("G27").Value = "c1, c3, c4, c9"
$ws.Range("G28").Value = "c1, c10"
$ws.Range("G29").Value = "c1, c8, c9"
$ws.Range("G30").Value = "c1, c2"

$gRange = $ws.Range("G25:I30")
$gRange.HorizontalAlignment = $xlLeft
$gRange.VerticalAlignment = $xlCenter

$ws.Range("G25:I25").Merge()
$ws.Range("G26:I26").Merge()
$ws.Range("G27:I27").Merge()
$ws.Range("G28:I28").Merge()
$ws.Range("G29:I29").Merge()
$ws.Range("G30:I30").Merge()

# "Detection Accuracy" column - stored as text fractions (e.g. "4/10")
$jRange = $ws.Range("J25:J30")
$jRange.NumberFormat = "@"
$ws.Range("J25").Value = "4/10"
$ws.Range("J26").Value = "5/10"
$ws.Range("J27").Value = "4/10"
$ws.Range("J28").Value = "2/10"
$ws.Range("J29").Value = "3/10"
$ws.Range("J30").Value = "2/10"
$jRange.HorizontalAlignment = $xlCenter
$jRange.VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 3. Chart title
# ---------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Average Detection Accuracy"

# ---------------------------------------------------------------------
# 4. Page setup - portrait orientation
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 5. Selection (matches the saved selection in the target workbook)
# ---------------------------------------------------------------------
$ws.Range("L18").Select()
